$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns (G:H) to hold the new error message columns ---
$ws.Range("G1:H1").EntireColumn.Insert()

# --- Remember the data that needs to move before clearing/rebuilding the grid ---

# Header labels: after the insert, D=Expected_lot_options, E=Population_name,
# F=Expected_File_names, I=Files_to_upload (old G, shifted right by the insert)
$hdrLotOptions = $ws.Cells.Item(1, 4).Value2
$hdrPopulation = $ws.Cells.Item(1, 5).Value2
$hdrExpectedFiles = $ws.Cells.Item(1, 6).Value2
$hdrFilesUpload = $ws.Cells.Item(1, 9).Value2

# "Expected_lot_options" list values, currently in column D, rows 2-9
$lotOptions = @()
for ($r = 2; $r -le 9; $r++) {
    $lotOptions += $ws.Cells.Item($r, 4).Value2
}

# The pop2 record currently lives in row 3:
#   A:C -> Name/LOT_name/Expected_ui_elements, D -> old list value (not needed here),
#   E:F -> Population_name/Expected_File_names, I -> Files_to_upload (shifted right by
#   the column insert above)
$pop2Name         = $ws.Cells.Item(3, 1).Value2
$pop2LotName      = $ws.Cells.Item(3, 2).Value2
$pop2UiElements   = $ws.Cells.Item(3, 3).Value2
$pop2Population   = $ws.Cells.Item(3, 5).Value2
$pop2ExpectedFile = $ws.Cells.Item(3, 6).Value2
$pop2FilesUpload  = $ws.Cells.Item(3, 9).Value2

# --- Clear the old layout so it can be rebuilt in its new shape ---
$ws.Range("A1:I9").ClearContents()

# --- Header row ---
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "LOT_name"
$ws.Range("C1").Value = "Expected_ui_elements"
$ws.Range("D1").Value = $hdrPopulation
$ws.Range("E1").Value = $hdrExpectedFiles
$ws.Range("F1").Value = $hdrFilesUpload
$ws.Range("G1").Value = "error_msg_col1"
$ws.Range("H1").Value = "error_msg_col2"
$ws.Range("I1").Value = $hdrLotOptions

# --- pop1 record (row 2) ---
$ws.Range("A2").Value = "pop1"
$ws.Range("B2").Value = "Automation_1"
$ws.Range("C2").Value = "Manage Line of Therapy"
$ws.Range("D2").Value = "Test - Test - 10/30/2020"
$ws.Range("E2").Value = "Test dataset - LineofTherapy_1.xlsx"
$ws.Range("F2").Value = "\Testdata\Templates\LineOfTherapy\Staging_Env\Test dataset - LineofTherapy_1.xlsx"
$ws.Range("G2").Value = 59
$ws.Range("H2").Value = "Line of Therapy 'Automation_1' is not supported"

# --- pop2 record moves down to row 4 ---
$ws.Range("A4").Value = $pop2Name
$ws.Range("B4").Value = $pop2LotName
$ws.Range("C4").Value = $pop2UiElements
$ws.Range("D4").Value = $pop2Population
$ws.Range("E4").Value = $pop2ExpectedFile
$ws.Range("F4").Value = $pop2FilesUpload
$ws.Range("G4").Value = 59
$ws.Range("H4").Value = "Line of Therapy 'Automation_2' is not supported"

# --- "Expected_lot_options" list moves to column I (same rows as before) ---
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 9).Value = $lotOptions[$r - 2]
}

# --- Column widths to match the resulting layout ---
$ws.Range("D1").ColumnWidth = 20.77734375
$ws.Range("E1").ColumnWidth = 28
$ws.Range("F1").ColumnWidth = 70.109375
$ws.Range("I1").ColumnWidth = 19

# --- Match the saved selection state ---
$ws.Range("A4:H4").Select()
